$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph that immediately follows
#    the title heading (it is being dropped entirely).
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
if ($metaPara.Range.Text -like "Meta description*") {
    $metaPara.Range.Delete()
}

# ------------------------------------------------------------------
# 2. Before the final "Prompt: ..." paragraph, insert a new paragraph
#    containing the bold title text, and replace the prompt text with
#    the new meta-description-style sentence (still italic).
# ------------------------------------------------------------------
$boldText = "Play Colossal Gems Free - Low Volatility With High RTP"

$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$insertPos = $lastPara.Range.Start

# Insert just the new bold paragraph's content; because InsertXML merges
# the final paragraph of the fragment with the paragraph at the
# insertion point, this effectively prepends our text onto the existing
# "Prompt: ..." paragraph instead of creating a break.
$insertRange = $d.Range($insertPos, $insertPos)
$xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>' + $boldText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertRange.InsertXML($xmlFrag)

# Now split the merged paragraph back into two: the new bold paragraph,
# and the original (still untouched) prompt paragraph.
$splitPos = $insertPos + $boldText.Length
$splitRange = $d.Range($splitPos, $splitPos)
$splitRange.InsertParagraphAfter()

# ------------------------------------------------------------------
# 3. Replace the old "Prompt: ..." text with the new description text,
#    keeping the existing italic run formatting and leading empty run.
# ------------------------------------------------------------------
$oldPromptText = "Prompt: Create a feature image for " + [char]34 + "Colossal Gems" + [char]34 + " in a cartoon style featuring a happy Maya warrior with glasses. DALLE, create a vibrant and colorful feature image for " + [char]34 + "Colossal Gems" + [char]34 + " that captures the essence of this joyful gem-themed slot game. Make sure to include a happy Maya warrior with glasses in the image to highlight its adventurous and playful nature. Use bright colors like pink, orange, and purple to give the image an eye-catching and dynamic look. Incorporate giant gem symbols of mega sizes, such as 2x2, 3x3, and 4x4, in the background to showcase the exciting bonus features of the game. Make it fun and inviting to encourage players to try their luck with this low volatility slot."
$newPromptText = "Try your luck on Colossal Gems, a gem-themed adventure with an Autospin function. This free slot game features low volatility and high RTP."

[void]$d.Content.Find.Execute(
    $oldPromptText,
    $true, $false, $false, $false, $false, $true, 1, $false,
    $newPromptText,
    2
)

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
